$d = $word.ActiveDocument

# Locate the target paragraph (the final "4. c. El ArrayList..." paragraph) robustly via Find.
$rng = $d.Content
$found = $rng.Find.Execute("4. c. El ArrayList", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate target paragraph"
}
$targetStart = $rng.Start

$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -eq $targetStart) {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq -1) {
    throw "Could not resolve paragraph index"
}

$emptyPara = $d.Paragraphs.Item($targetIndex - 1)
$oldPara = $d.Paragraphs.Item($targetIndex)

# New OOXML block: point 3 (a/b + two Utility Class code examples) followed by
# the corrected point 4.c. paragraph (with the bookmark now anchored right
# after the final code paragraph's closing brace run).
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>3.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:br/></w:r><w:r><w:t>a. Las U</w:t></w:r><w:r><w:t>ti</w:t></w:r><w:r><w:t>li</w:t></w:r><w:r><w:t xml:space="preserve">ly </w:t></w:r><w:r><w:t>C</w:t></w:r><w:r><w:t xml:space="preserve">lass son clases </w:t></w:r><w:r><w:t xml:space="preserve">que se usan muy a menudo, estas se conforman siempre de atributos estáticos, y </w:t></w:r><w:r><w:t>deben ser públicas para poder ser accedidas desde cualquier parte</w:t></w:r><w:r><w:t xml:space="preserve"> del proyecto</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> también deben tener al inicio la palabra reservada Final,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>para que no puedan ser editada</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t xml:space="preserve">y por ultimo </w:t></w:r><w:r><w:t>sus métodos también tienen que ser estáticos</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t>b.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Ejemplo de una Utility Class:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Usamos</w:t></w:r><w:r><w:t xml:space="preserve"> la clase Math </w:t></w:r><w:r><w:t>para hallar la raíz de un numero ingresado</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Import java.util.Math;</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Import java.util.Scanner; </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Public class ClassMath{</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Scanner sc= new Scanner(System.in); </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>System.out.println(“Por</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> favor digite un </w:t></w:r><w:r><w:t>número</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, al cual le sacaremos la </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>raíz</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> cuadrada: </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">”); </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Int </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>num</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>=sc. nextInt();</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Int aux;</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">     </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>If</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(num =&gt; 0){</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">  </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">      </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> aux = sqrt(num);</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">      </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>System.out.println(“La raiz cuadrada de “+num+” es igual a: “+aux);</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br/><w:t>}</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Ahora vamos hacer un tipo de Utility Class, el cual nos halla </w:t></w:r><w:r><w:t>fácilmente</w:t></w:r><w:r><w:t xml:space="preserve"> el </w:t></w:r><w:r><w:t>área</w:t></w:r><w:r><w:t xml:space="preserve"> de un Triangulo:</w:t></w:r><w:r><w:br/></w:r><w:r><w:br/><w:t>import java.util.Scanner;  </w:t></w:r></w:p><w:p><w:r><w:t>public class </w:t></w:r><w:r><w:t>TestUtility</w:t></w:r><w:r><w:t> {  </w:t></w:r></w:p><w:p><w:r><w:t>     public static void main(String[] args) {  </w:t></w:r></w:p><w:p><w:r><w:t>        Scanner sc = new Scanner(System.in);  </w:t></w:r></w:p><w:p><w:r><w:t>        double b,h;  </w:t></w:r></w:p><w:p><w:r><w:t>        System.out.println(</w:t></w:r><w:r><w:t>"Ingrese</w:t></w:r><w:r><w:t> </w:t></w:r><w:r><w:t xml:space="preserve">la </w:t></w:r><w:r><w:t>base</w:t></w:r><w:r><w:t xml:space="preserve"> del triangulo</w:t></w:r><w:r><w:t>");  </w:t></w:r></w:p><w:p><w:r><w:t>        b=sc.nextDouble();  </w:t></w:r></w:p><w:p><w:r><w:t>        System.out.println("Ingresa </w:t></w:r><w:r><w:t xml:space="preserve">la </w:t></w:r><w:r><w:t>altura</w:t></w:r><w:r><w:t xml:space="preserve"> del triangulo</w:t></w:r><w:r><w:t>");  </w:t></w:r></w:p><w:p><w:r><w:t>        h=sc.nextDouble();  </w:t></w:r></w:p><w:p><w:r><w:t>        double area;  </w:t></w:r></w:p><w:p><w:r><w:t>        area=b*h/2;  </w:t></w:r></w:p><w:p><w:r><w:t>       System.out.print</w:t></w:r><w:r><w:t>ln</w:t></w:r><w:r><w:t>(</w:t></w:r><w:r><w:t>“Este es el área del triangulo: “+</w:t></w:r><w:r><w:t>area);  </w:t></w:r></w:p><w:p><w:r><w:t>     }  </w:t></w:r></w:p><w:p><w:r><w:t>}  </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:r><w:t>4. c. El ArrayList debe de ser del tipo “nombre de la clase padre” la cual en este caso es la abstracta, ya que es el ítem que tienen todos en común, porque son “parte a la vez del padre”</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>
'@

# Insert the whole block in place of the empty paragraph; this pushes the old
# "4. c. ..." paragraph further down, unchanged.
$emptyPara.Range.InsertXML($xml)

# The old "4. c. ..." paragraph (with its stray bookmark) is now redundant --
# the new block's last paragraph already carries the corrected text and the
# relocated bookmark. Remove its contents; since it is the final paragraph in
# the body, Word collapses the now-empty trailing mark away automatically.
$oldPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$prevPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$delRange = $d.Range($prevPara.Range.End, $oldPara.Range.End)
$delRange.Delete()

Write-Output "done"
